$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("avnnet")
$ws2 = $wb.Worksheets.Item("bagging")
$ws = $ws2

# --- AVNNET rows 22-41 (duplicate of avnnet sheet data; reuses existing shared strings) ---
$ws.Cells.Item(22,1).NumberFormat = "@"
$ws.Cells.Item(22,1).Value = "0.1067646"
$ws.Cells.Item(22,1).Style = "Normal"
$ws.Cells.Item(22,2).NumberFormat = "@"
$ws.Cells.Item(22,2).Value = "0.8960432"
$ws.Cells.Item(22,2).Style = "Normal"
$ws.Cells.Item(23,1).NumberFormat = "@"
$ws.Cells.Item(23,1).Value = "0.1062521"
$ws.Cells.Item(23,1).Style = "Normal"
$ws.Cells.Item(23,2).NumberFormat = "@"
$ws.Cells.Item(23,2).Value = "0.8974809"
$ws.Cells.Item(23,2).Style = "Normal"
$ws.Cells.Item(24,1).NumberFormat = "@"
$ws.Cells.Item(24,1).Value = "0.1055688"
$ws.Cells.Item(24,1).Style = "Normal"
$ws.Cells.Item(24,2).NumberFormat = "@"
$ws.Cells.Item(24,2).Value = "0.8947908"
$ws.Cells.Item(24,2).Style = "Normal"
$ws.Cells.Item(25,1).NumberFormat = "@"
$ws.Cells.Item(25,1).Value = "0.1064230"
$ws.Cells.Item(25,1).Style = "Normal"
$ws.Cells.Item(25,2).NumberFormat = "@"
$ws.Cells.Item(25,2).Value = "0.8982033"
$ws.Cells.Item(25,2).Style = "Normal"
$ws.Cells.Item(26,1).NumberFormat = "@"
$ws.Cells.Item(26,1).Value = "0.1045439"
$ws.Cells.Item(26,1).Style = "Normal"
$ws.Cells.Item(26,2).NumberFormat = "@"
$ws.Cells.Item(26,2).Value = "0.8964518"
$ws.Cells.Item(26,2).Style = "Normal"
$ws.Cells.Item(27,1).NumberFormat = "@"
$ws.Cells.Item(27,1).Value = "0.1062521"
$ws.Cells.Item(27,1).Style = "Normal"
$ws.Cells.Item(27,2).NumberFormat = "@"
$ws.Cells.Item(27,2).Value = "0.8965810"
$ws.Cells.Item(27,2).Style = "Normal"
$ws.Cells.Item(28,1).NumberFormat = "@"
$ws.Cells.Item(28,1).Value = "0.1059105"
$ws.Cells.Item(28,1).Style = "Normal"
$ws.Cells.Item(28,2).NumberFormat = "@"
$ws.Cells.Item(28,2).Value = "0.8945043"
$ws.Cells.Item(28,2).Style = "Normal"
$ws.Cells.Item(29,1).NumberFormat = "@"
$ws.Cells.Item(29,1).Value = "0.1065938"
$ws.Cells.Item(29,1).Style = "Normal"
$ws.Cells.Item(29,2).NumberFormat = "@"
$ws.Cells.Item(29,2).Value = "0.8990411"
$ws.Cells.Item(29,2).Style = "Normal"
$ws.Cells.Item(30,1).NumberFormat = "@"
$ws.Cells.Item(30,1).Value = "0.1033481"
$ws.Cells.Item(30,1).Style = "Normal"
$ws.Cells.Item(30,2).NumberFormat = "@"
$ws.Cells.Item(30,2).Value = "0.8969116"
$ws.Cells.Item(30,2).Style = "Normal"
$ws.Cells.Item(31,1).NumberFormat = "@"
$ws.Cells.Item(31,1).Value = "0.1055688"
$ws.Cells.Item(31,1).Style = "Normal"
$ws.Cells.Item(31,2).NumberFormat = "@"
$ws.Cells.Item(31,2).Value = "0.8957565"
$ws.Cells.Item(31,2).Style = "Normal"
$ws.Cells.Item(32,1).NumberFormat = "@"
$ws.Cells.Item(32,1).Value = "0.1088145"
$ws.Cells.Item(32,1).Style = "Normal"
$ws.Cells.Item(32,2).NumberFormat = "@"
$ws.Cells.Item(32,2).Value = "0.8980774"
$ws.Cells.Item(32,2).Style = "Normal"
$ws.Cells.Item(33,1).NumberFormat = "@"
$ws.Cells.Item(33,1).Value = "0.1141100"
$ws.Cells.Item(33,1).Style = "Normal"
$ws.Cells.Item(33,2).NumberFormat = "@"
$ws.Cells.Item(33,2).Value = "0.8969203"
$ws.Cells.Item(33,2).Style = "Normal"
$ws.Cells.Item(34,1).NumberFormat = "@"
$ws.Cells.Item(34,1).Value = "0.1110352"
$ws.Cells.Item(34,1).Style = "Normal"
$ws.Cells.Item(34,2).NumberFormat = "@"
$ws.Cells.Item(34,2).Value = "0.9028122"
$ws.Cells.Item(34,2).Style = "Normal"
$ws.Cells.Item(35,1).NumberFormat = "@"
$ws.Cells.Item(35,1).Value = "0.1100102"
$ws.Cells.Item(35,1).Style = "Normal"
$ws.Cells.Item(35,2).NumberFormat = "@"
$ws.Cells.Item(35,2).Value = "0.8974087"
$ws.Cells.Item(35,2).Style = "Normal"
$ws.Cells.Item(36,1).NumberFormat = "@"
$ws.Cells.Item(36,1).Value = "0.1117185"
$ws.Cells.Item(36,1).Style = "Normal"
$ws.Cells.Item(36,2).NumberFormat = "@"
$ws.Cells.Item(36,2).Value = "0.8941729"
$ws.Cells.Item(36,2).Style = "Normal"
$ws.Cells.Item(37,1).NumberFormat = "@"
$ws.Cells.Item(37,1).Value = "0.1089853"
$ws.Cells.Item(37,1).Style = "Normal"
$ws.Cells.Item(37,2).NumberFormat = "@"
$ws.Cells.Item(37,2).Value = "0.8957019"
$ws.Cells.Item(37,2).Style = "Normal"
$ws.Cells.Item(38,1).NumberFormat = "@"
$ws.Cells.Item(38,1).Value = "0.1105227"
$ws.Cells.Item(38,1).Style = "Normal"
$ws.Cells.Item(38,2).NumberFormat = "@"
$ws.Cells.Item(38,2).Value = "0.8958447"
$ws.Cells.Item(38,2).Style = "Normal"
$ws.Cells.Item(39,1).NumberFormat = "@"
$ws.Cells.Item(39,1).Value = "0.1106935"
$ws.Cells.Item(39,1).Style = "Normal"
$ws.Cells.Item(39,2).NumberFormat = "@"
$ws.Cells.Item(39,2).Value = "0.8996416"
$ws.Cells.Item(39,2).Style = "Normal"
$ws.Cells.Item(40,1).NumberFormat = "@"
$ws.Cells.Item(40,1).Value = "0.1153058"
$ws.Cells.Item(40,1).Style = "Normal"
$ws.Cells.Item(40,2).NumberFormat = "@"
$ws.Cells.Item(40,2).Value = "0.8968740"
$ws.Cells.Item(40,2).Style = "Normal"
$ws.Cells.Item(41,1).NumberFormat = "@"
$ws.Cells.Item(41,1).Value = "0.1118893"
$ws.Cells.Item(41,1).Style = "Normal"
$ws.Cells.Item(41,2).NumberFormat = "@"
$ws.Cells.Item(41,2).Value = "0.8967188"
$ws.Cells.Item(41,2).Style = "Normal"
$ws.Cells.Item(22,3).Value = "AVNNET MODELO 1"
$ws.Cells.Item(23,3).Value = "AVNNET MODELO 1"
$ws.Cells.Item(24,3).Value = "AVNNET MODELO 1"
$ws.Cells.Item(25,3).Value = "AVNNET MODELO 1"
$ws.Cells.Item(26,3).Value = "AVNNET MODELO 1"
$ws.Cells.Item(27,3).Value = "AVNNET MODELO 1"
$ws.Cells.Item(28,3).Value = "AVNNET MODELO 1"
$ws.Cells.Item(29,3).Value = "AVNNET MODELO 1"
$ws.Cells.Item(30,3).Value = "AVNNET MODELO 1"
$ws.Cells.Item(31,3).Value = "AVNNET MODELO 1"
$ws.Cells.Item(32,3).Value = "AVNNET MODELO 2"
$ws.Cells.Item(33,3).Value = "AVNNET MODELO 2"
$ws.Cells.Item(34,3).Value = "AVNNET MODELO 2"
$ws.Cells.Item(35,3).Value = "AVNNET MODELO 2"
$ws.Cells.Item(36,3).Value = "AVNNET MODELO 2"
$ws.Cells.Item(37,3).Value = "AVNNET MODELO 2"
$ws.Cells.Item(38,3).Value = "AVNNET MODELO 2"
$ws.Cells.Item(39,3).Value = "AVNNET MODELO 2"
$ws.Cells.Item(40,3).Value = "AVNNET MODELO 2"
$ws.Cells.Item(41,3).Value = "AVNNET MODELO 2"

# --- BAG. MODELO 2 rows 52-61 (written before MODELO 1 to match shared-string append order) ---
$ws.Cells.Item(52,1).NumberFormat = "@"
$ws.Cells.Item(52,1).Value = "0.10112743"
$ws.Cells.Item(52,1).Style = "Normal"
$ws.Cells.Item(52,2).NumberFormat = "@"
$ws.Cells.Item(52,2).Value = "0.9126257"
$ws.Cells.Item(52,2).Style = "Normal"
$ws.Cells.Item(53,1).NumberFormat = "@"
$ws.Cells.Item(53,1).Value = "0.09993167"
$ws.Cells.Item(53,1).Style = "Normal"
$ws.Cells.Item(53,2).NumberFormat = "@"
$ws.Cells.Item(53,2).Value = "0.9144836"
$ws.Cells.Item(53,2).Style = "Normal"
$ws.Cells.Item(54,1).NumberFormat = "@"
$ws.Cells.Item(54,1).Value = "0.09890673"
$ws.Cells.Item(54,1).Style = "Normal"
$ws.Cells.Item(54,2).NumberFormat = "@"
$ws.Cells.Item(54,2).Value = "0.9145415"
$ws.Cells.Item(54,2).Style = "Normal"
$ws.Cells.Item(55,1).NumberFormat = "@"
$ws.Cells.Item(55,1).Value = "0.09959002"
$ws.Cells.Item(55,1).Style = "Normal"
$ws.Cells.Item(55,2).NumberFormat = "@"
$ws.Cells.Item(55,2).Value = "0.9140330"
$ws.Cells.Item(55,2).Style = "Normal"
$ws.Cells.Item(56,1).NumberFormat = "@"
$ws.Cells.Item(56,1).Value = "0.10027332"
$ws.Cells.Item(56,1).Style = "Normal"
$ws.Cells.Item(56,2).NumberFormat = "@"
$ws.Cells.Item(56,2).Value = "0.9140948"
$ws.Cells.Item(56,2).Style = "Normal"
$ws.Cells.Item(57,1).NumberFormat = "@"
$ws.Cells.Item(57,1).Value = "0.09941920"
$ws.Cells.Item(57,1).Style = "Normal"
$ws.Cells.Item(57,2).NumberFormat = "@"
$ws.Cells.Item(57,2).Value = "0.9147203"
$ws.Cells.Item(57,2).Style = "Normal"
$ws.Cells.Item(58,1).NumberFormat = "@"
$ws.Cells.Item(58,1).Value = "0.10010249"
$ws.Cells.Item(58,1).Style = "Normal"
$ws.Cells.Item(58,2).NumberFormat = "@"
$ws.Cells.Item(58,2).Value = "0.9145659"
$ws.Cells.Item(58,2).Style = "Normal"
$ws.Cells.Item(59,1).NumberFormat = "@"
$ws.Cells.Item(59,1).Value = "0.09907755"
$ws.Cells.Item(59,1).Style = "Normal"
$ws.Cells.Item(59,2).NumberFormat = "@"
$ws.Cells.Item(59,2).Value = "0.9155725"
$ws.Cells.Item(59,2).Style = "Normal"
$ws.Cells.Item(60,1).NumberFormat = "@"
$ws.Cells.Item(60,1).Value = "0.10095661"
$ws.Cells.Item(60,1).Style = "Normal"
$ws.Cells.Item(60,2).NumberFormat = "@"
$ws.Cells.Item(60,2).Value = "0.9155204"
$ws.Cells.Item(60,2).Style = "Normal"
$ws.Cells.Item(61,1).NumberFormat = "@"
$ws.Cells.Item(61,1).Value = "0.10112743"
$ws.Cells.Item(61,1).Style = "Normal"
$ws.Cells.Item(61,2).NumberFormat = "@"
$ws.Cells.Item(61,2).Value = "0.9127317"
$ws.Cells.Item(61,2).Style = "Normal"
$ws.Cells.Item(52,3).Value = "BAG. MODELO 2"
$ws.Cells.Item(53,3).Value = "BAG. MODELO 2"
$ws.Cells.Item(54,3).Value = "BAG. MODELO 2"
$ws.Cells.Item(55,3).Value = "BAG. MODELO 2"
$ws.Cells.Item(56,3).Value = "BAG. MODELO 2"
$ws.Cells.Item(57,3).Value = "BAG. MODELO 2"
$ws.Cells.Item(58,3).Value = "BAG. MODELO 2"
$ws.Cells.Item(59,3).Value = "BAG. MODELO 2"
$ws.Cells.Item(60,3).Value = "BAG. MODELO 2"
$ws.Cells.Item(61,3).Value = "BAG. MODELO 2"

# --- BAG. MODELO 1 rows 42-51 ---
$ws.Cells.Item(42,1).NumberFormat = "@"
$ws.Cells.Item(42,1).Value = "0.10181073"
$ws.Cells.Item(42,1).Style = "Normal"
$ws.Cells.Item(42,2).NumberFormat = "@"
$ws.Cells.Item(42,2).Value = "0.9124968"
$ws.Cells.Item(42,2).Style = "Normal"
$ws.Cells.Item(43,1).NumberFormat = "@"
$ws.Cells.Item(43,1).Value = "0.10317731"
$ws.Cells.Item(43,1).Style = "Normal"
$ws.Cells.Item(43,2).NumberFormat = "@"
$ws.Cells.Item(43,2).Value = "0.9144305"
$ws.Cells.Item(43,2).Style = "Normal"
$ws.Cells.Item(44,1).NumberFormat = "@"
$ws.Cells.Item(44,1).Value = "0.10078579"
$ws.Cells.Item(44,1).Style = "Normal"
$ws.Cells.Item(44,2).NumberFormat = "@"
$ws.Cells.Item(44,2).Value = "0.9143850"
$ws.Cells.Item(44,2).Style = "Normal"
$ws.Cells.Item(45,1).NumberFormat = "@"
$ws.Cells.Item(45,1).Value = "0.10181073"
$ws.Cells.Item(45,1).Style = "Normal"
$ws.Cells.Item(45,2).NumberFormat = "@"
$ws.Cells.Item(45,2).Value = "0.9138405"
$ws.Cells.Item(45,2).Style = "Normal"
$ws.Cells.Item(46,1).NumberFormat = "@"
$ws.Cells.Item(46,1).Value = "0.10198155"
$ws.Cells.Item(46,1).Style = "Normal"
$ws.Cells.Item(46,2).NumberFormat = "@"
$ws.Cells.Item(46,2).Value = "0.9144000"
$ws.Cells.Item(46,2).Style = "Normal"
$ws.Cells.Item(47,1).NumberFormat = "@"
$ws.Cells.Item(47,1).Value = "0.10061496"
$ws.Cells.Item(47,1).Style = "Normal"
$ws.Cells.Item(47,2).NumberFormat = "@"
$ws.Cells.Item(47,2).Value = "0.9151003"
$ws.Cells.Item(47,2).Style = "Normal"
$ws.Cells.Item(48,1).NumberFormat = "@"
$ws.Cells.Item(48,1).Value = "0.10249402"
$ws.Cells.Item(48,1).Style = "Normal"
$ws.Cells.Item(48,2).NumberFormat = "@"
$ws.Cells.Item(48,2).Value = "0.9151371"
$ws.Cells.Item(48,2).Style = "Normal"
$ws.Cells.Item(49,1).NumberFormat = "@"
$ws.Cells.Item(49,1).Value = "0.10163990"
$ws.Cells.Item(49,1).Style = "Normal"
$ws.Cells.Item(49,2).NumberFormat = "@"
$ws.Cells.Item(49,2).Value = "0.9150177"
$ws.Cells.Item(49,2).Style = "Normal"
$ws.Cells.Item(50,1).NumberFormat = "@"
$ws.Cells.Item(50,1).Value = "0.10129826"
$ws.Cells.Item(50,1).Style = "Normal"
$ws.Cells.Item(50,2).NumberFormat = "@"
$ws.Cells.Item(50,2).Value = "0.9149272"
$ws.Cells.Item(50,2).Style = "Normal"
$ws.Cells.Item(51,1).NumberFormat = "@"
$ws.Cells.Item(51,1).Value = "0.10181073"
$ws.Cells.Item(51,1).Style = "Normal"
$ws.Cells.Item(51,2).NumberFormat = "@"
$ws.Cells.Item(51,2).Value = "0.9136212"
$ws.Cells.Item(51,2).Style = "Normal"
$ws.Cells.Item(42,3).Value = "BAG. MODELO 1"
$ws.Cells.Item(43,3).Value = "BAG. MODELO 1"
$ws.Cells.Item(44,3).Value = "BAG. MODELO 1"
$ws.Cells.Item(45,3).Value = "BAG. MODELO 1"
$ws.Cells.Item(46,3).Value = "BAG. MODELO 1"
$ws.Cells.Item(47,3).Value = "BAG. MODELO 1"
$ws.Cells.Item(48,3).Value = "BAG. MODELO 1"
$ws.Cells.Item(49,3).Value = "BAG. MODELO 1"
$ws.Cells.Item(50,3).Value = "BAG. MODELO 1"
$ws.Cells.Item(51,3).Value = "BAG. MODELO 1"

# --- sheet view / selection changes ---
$ws1.Activate()
$excel.ActiveWindow.Zoom = 92
$ws1.Range("C45").Select() | Out-Null

$ws2.Activate()
$ws2.Range("A42:C61").Select() | Out-Null

Write-Host "Edit applied"
